$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column D ("http://dbpedia.org/ontology/birthPlace") entirely -
# it (and every per-row blank cell under it) disappears in the target sheet.
$ws.Columns.Item(4).Delete()

# Append the new resource rows (43-71) with their birth/death date columns.
$newRows = @(
    @("http://dbpedia.org/resource/Francis_Freeling", $null, $null),
    @("http://dbpedia.org/resource/Francis_Champneys", "1848-03-25 ", "1930-07-30 "),
    @("http://dbpedia.org/resource/Ferdinand_Dalberg-Acton", $null, $null),
    @("http://dbpedia.org/resource/Ernest_Musgrave_Harvey", $null, $null),
    @("http://dbpedia.org/resource/Ernest_Craig", $null, $null),
    @("http://dbpedia.org/resource/Erik_Ohlson", $null, $null),
    @("http://dbpedia.org/resource/Edward_des_Bouverie", $null, $null),
    @("http://dbpedia.org/resource/Edward_Manningham-Buller", $null, $null),
    @("http://dbpedia.org/resource/Edward_Mackay_Edgar", $null, $null),
    @("http://dbpedia.org/resource/Edward_Irby", $null, $null),
    @("http://dbpedia.org/resource/Edmund_Findlay", "1902-06-14 ", "1962-09-06 "),
    @("http://dbpedia.org/resource/Cuthbert_Ackroyd", $null, $null),
    @("http://dbpedia.org/resource/Currimbhoy_Ebrahim", $null, $null),
    @("http://dbpedia.org/resource/Cory_Cory-Wright", $null, $null),
    @("http://dbpedia.org/resource/Charles_Renshaw", $null, $null),
    @("http://dbpedia.org/resource/Charles_Jessel", $null, $null),
    @("http://dbpedia.org/resource/Charles_Gladstone", $null, $null),
    @("http://dbpedia.org/resource/Charles_Cornwallis_Lloyd", $null, $null),
    @("http://dbpedia.org/resource/Brograve_Beauchamp", $null, $null),
    @("http://dbpedia.org/resource/Bernard_Waley-Cohen", $null, $null),
    @("http://dbpedia.org/resource/August_Cayzer", $null, $null),
    @("http://dbpedia.org/resource/Arthur_Cory-Wright", $null, $null),
    @("http://dbpedia.org/resource/Archibald_Birkmyre", $null, $null),
    @("http://dbpedia.org/resource/Antony_Guy_Acland", $null, $null),
    @("http://dbpedia.org/resource/Anthony_Aucher", $null, $null),
    @("http://dbpedia.org/resource/Andrew_Lauder", $null, $null),
    @("http://dbpedia.org/resource/Andrew_Armstrong", $null, $null),
    @("http://dbpedia.org/resource/Alexander_Erskine-Hill", $null, $null),
    @("http://dbpedia.org/resource/Abraham_Janssen", $null, $null)
)

$row = 43
foreach ($r in $newRows) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    if ($r[1] -ne $null) { $ws.Cells.Item($row, 2).Value = $r[1] }
    if ($r[2] -ne $null) { $ws.Cells.Item($row, 3).Value = $r[2] }
    $row = $row + 1
}
